# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx data refresh described in the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''24.674.53'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.38%  '

# Row 3
$ws.Range("D3").Value = '''1.695.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.32%  '

# Row 4
$ws.Range("E4").Value = '  +0.20%  '

# Row 5
$ws.Range("D5").Value = '''317.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.85%  '

# Row 6
$ws.Range("E6").Value = '  +0.16%  '

# Row 7
$ws.Range("D7").Value = '''0.3963'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.04%  '

# Row 8
$ws.Range("D8").Value = '''0.4083'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.88%  '

# Row 9
$ws.Range("D9").Value = '''1.498'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.29%  '

# Row 10
$ws.Range("D10").Value = '''1.001'
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = '''51.47'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.02%  '

# Row 12
$ws.Range("D12").Value = '''0.08951'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.18%  '

# Row 13
$ws.Range("D13").Value = '''7.205'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.34%  '

# Row 14
$ws.Range("D14").Value = '''23.57'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.38%  '

# Row 15
$ws.Range("D15").Value = '''8.207'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +11.94%  '

# Row 16
$ws.Range("D16").Value = '''0.00001333'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.80%  '

# Row 17
$ws.Range("D17").Value = '''1.696.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.64%  '

# Row 18
$ws.Range("D18").Value = '''100.17'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.13%  '

# Row 19
$ws.Range("D19").Value = '''0.07007'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.06%  '

# Row 20
$ws.Range("D20").Value = '''19.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.03%  '

# Row 21
$ws.Range("D21").Value = '''7.045'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.23%  '

# Row 22
$ws.Range("D22").Value = '''1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.23%  '

# Row 23
$ws.Range("D23").Value = '''14.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.93%  '

# Row 24
$ws.Range("D24").Value = '''24.679.69'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.44%  '

# Row 25
$ws.Range("D25").Value = '''3.143'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.66%  '

# Row 26
$ws.Range("D26").Value = '''2.343'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.16%  '

# Row 27
$ws.Range("D27").Value = '''22.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.56%  '

# Row 28
$ws.Range("D28").Value = '''163.46'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.74%  '

# Row 29
$ws.Range("D29").Value = '''137.50'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.42%  '

# Row 30
$ws.Range("D30").Value = '''5.170'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.80%  '

# Row 31
$ws.Range("D31").Value = '''7.430'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.64%  '

# Row 32
$ws.Range("D32").Value = '''1.879.13'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.39%  '

# Row 33
$ws.Range("D33").Value = '''1.075'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.13%  '

# Row 34
$ws.Range("D34").Value = '''0.08613'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.06%  '

# Row 35
$ws.Range("D35").Value = '''7.151'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.06%  '

# Row 36
$ws.Range("D36").Value = '''11.46'
$ws.Range("D36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = '''0.2749'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.71%  '

# Row 38
$ws.Range("E38").Value = '  +0.88%  '

# Row 39
$ws.Range("D39").Value = '''14.47'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.05%  '

# Row 40
$ws.Range("D40").Value = '''0.09226'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.48%  '

# Row 41
$ws.Range("D41").Value = '''0.02731'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.41%  '

# Row 42
$ws.Range("D42").Value = '''1.479'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.80%  '

# Row 43
$ws.Range("D43").Value = '''0.7691'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.25%  '

# Row 44
$ws.Range("D44").Value = '''16.04'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.55%  '

# Row 45
$ws.Range("D45").Value = '''2.629'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.69%  '

# Row 46
$ws.Range("D46").Value = '''0.7186'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.68%  '

# Row 47
$ws.Range("D47").Value = '''4.229'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.08%  '

# Row 48
$ws.Range("D48").Value = '''1.001'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.21%  '

# Row 49
$ws.Range("B49").Value = 'Flow'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D49").Value = '''1.327'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.63%  '

# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '''140.56'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.14%  '

# Row 51
$ws.Range("D51").Value = '''0.07982'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.61%  '

